# Refresh the scraped cryptocurrency price/volume snapshot (GitHub Actions cron update).
# Only the Price (D) / Volume(1h) (E) columns change value for most rows; rows 47-48
# (Stellar / ApeXProtocol) also swap rank position, so their Coin/Link cells change too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.808.96'
$ws.Range('E2').Value = '  +4.12%  '
$ws.Range('D3').Value = '3.437.75'
$ws.Range('E3').Value = '  +3.56%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''578.80'
$ws.Range('E5').Value = '  +4.59%  '
$ws.Range('D6').Value = '''184.05'
$ws.Range('E6').Value = '  +6.10%  '
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('D8').Value = '3.430.20'
$ws.Range('E8').Value = '  +3.57%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').Value = '''56.15'
$ws.Range('E12').Value = '  +4.86%  '
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '''9.41'
$ws.Range('E14').Value = '  +4.22%  '
$ws.Range('D15').Value = '3.991.46'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('E16').Value = '  +3.27%  '
$ws.Range('D17').Value = '3.444.23'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = '66.697.16'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('E20').Value = '  +3.49%  '
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').Value = '''483.50'
$ws.Range('E22').Value = '  +6.87%  '
$ws.Range('D23').Value = '''16.95'
$ws.Range('E23').Value = '  +23.44%  '
$ws.Range('D24').Value = '''5.13'
$ws.Range('E24').Value = '  +3.25%  '
$ws.Range('E25').Value = '  +7.61%  '
$ws.Range('D26').Value = '''89.63'
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('D27').Value = '''11.07'
$ws.Range('E27').Value = '  +3.82%  '
$ws.Range('D28').Value = '''2.94'
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('E29').Value = '  +6.87%  '
$ws.Range('D30').Value = '''31.28'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('E31').Value = '  +8.80%  '
$ws.Range('D32').Value = '''64.40'
$ws.Range('E32').Value = '  +6.30%  '
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('D34').Value = '''591.38'
$ws.Range('E34').Value = '  +4.49%  '
$ws.Range('E35').Value = '  +4.72%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +5.58%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = '''36.40'
$ws.Range('E39').Value = '  +3.65%  '
$ws.Range('E40').Value = '  +5.53%  '
$ws.Range('E41').Value = '  +4.72%  '
$ws.Range('D42').Value = '3.188.58'
$ws.Range('E42').Value = '  +4.38%  '
$ws.Range('E43').Value = '  +5.37%  '
$ws.Range('E44').Value = '  +3.98%  '
$ws.Range('D45').Value = '''2.54'
$ws.Range('E45').Value = '  +5.00%  '
$ws.Range('D46').Value = '''2.77'
$ws.Range('E46').Value = '  +22.08%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''3.22'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '''0.135'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').Value = '''8.71'
$ws.Range('E50').Value = '  +6.93%  '
$ws.Range('D51').Value = '''139.79'
$ws.Range('E51').Value = '  -2.08%  '
